# Add two more "run" result rows (rows 6 and 7) to the Statistic sheet.
#
# Like the existing row 4 (configuration #2), these new rows only have
# results for the first run block ("Прогін 0", columns A:J) and the
# overall success-rate column (AC); the remaining run blocks (K:AB,
# "Прогін 1" and "Прогін 2") stay blank/empty for these configurations.
#
# We seed each new row by copying the blank-cell layout from row 4 (so
# the otherwise-empty K:AB cells are materialized exactly like the rest
# of the sheet does for single-run configurations), then overwrite the
# cells that actually carry data for the new configurations.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @{ Row = 6; A = 4; B = 196700; C = 263; D = 5; E = 1; F = 4; G = 0.019011406844106463; H = 1; I = 0.4; J = 0.9809885931558935; AC = 0 },
    @{ Row = 7; A = 5; B = 196700; C = 263; D = 5; E = 1; F = 4; G = 0.019011406844106463; H = 1; I = 0.4; J = 0.9809885931558935; AC = 0 }
)

foreach ($rowData in $newRows) {
    $r = $rowData.Row

    # Bring over the blank A:AC cell layout used by row 4.
    $ws.Range("A4:AC4").Copy($ws.Range("A$($r):AC$($r)"))

    # Fill in this configuration's actual values.
    $ws.Range("A$r").Value = $rowData.A
    $ws.Range("B$r").Value = $rowData.B
    $ws.Range("C$r").Value = $rowData.C
    $ws.Range("D$r").Value = $rowData.D
    $ws.Range("E$r").Value = $rowData.E
    $ws.Range("F$r").Value = $rowData.F
    $ws.Range("G$r").Value = $rowData.G
    $ws.Range("H$r").Value = $rowData.H
    $ws.Range("I$r").Value = $rowData.I
    $ws.Range("J$r").Value = $rowData.J
    $ws.Range("AC$r").Value = $rowData.AC
}
